# ============================================================================
# Applies the "additional scraping" edit:
#   - Adds a new "Player Info" sheet (becomes the first sheet).
#   - Keeps "ODI Batting" and "ODI Bowling" (now 2nd / 3rd), with:
#       * MATCH_CARD_LINK column renamed to MATCH_CODE
#       * the URL values replaced by just the trailing MatchCode number
#       * leftover empty INNING_NUMBER cells (ODI Batting col B) cleared out
#   - Adds a new "ODI Batting Extra" sheet (becomes the 4th / last sheet).
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# 1) ODI Batting: MATCH_CARD_LINK -> MATCH_CODE (header + values), and
#    drop the stray empty INNING_NUMBER (column B) placeholder cells.
# ----------------------------------------------------------------------
$wsBatting = $wb.Worksheets.Item("ODI Batting")

$battingLastRow = 94
$wsBatting.Range("D2:D$battingLastRow").NumberFormat = "@"

for ($r = 2; $r -le $battingLastRow; $r++) {
    $linkCell = $wsBatting.Cells.Item($r, 4)
    $link = $linkCell.Value()
    $pieces = $link -split "MatchCode="
    $matchCode = $pieces[-1]
    $linkCell.Value = $matchCode

    $inningCell = $wsBatting.Cells.Item($r, 2)
    $inningVal = $inningCell.Value()
    if ($inningVal -eq $null -or $inningVal -eq "") {
        $inningCell.Value = ""
    }
}

$wsBatting.Range("D1").Value = "MATCH_CODE"

# ----------------------------------------------------------------------
# 2) ODI Bowling: MATCH_CARD_LINK -> MATCH_CODE (header + values).
# ----------------------------------------------------------------------
$wsBowling = $wb.Worksheets.Item("ODI Bowling")

$bowlingLastRow = 89
$wsBowling.Range("B2:B$bowlingLastRow").NumberFormat = "@"

for ($r = 2; $r -le $bowlingLastRow; $r++) {
    $linkCell = $wsBowling.Cells.Item($r, 2)
    $link = $linkCell.Value()
    $pieces = $link -split "MatchCode="
    $matchCode = $pieces[-1]
    $linkCell.Value = $matchCode
}

$wsBowling.Range("B1").Value = "MATCH_CODE"

# ----------------------------------------------------------------------
# 3) New "Player Info" sheet, placed before "ODI Batting" (first tab).
# ----------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Add()
$wsInfo.Name = "Player Info"

$wsInfo.Range("A1").Value = "ID"
$wsInfo.Range("B1").Value = "NAME"
$wsInfo.Range("C1").Value = "BATTING_HAND"
$wsInfo.Range("D1").Value = "BOWL_STYLE"

$wsInfo.Range("A1:D1").Style = $wsBatting.Range("A1").Style

$wsInfo.Range("A2").NumberFormat = "@"
$wsInfo.Range("A2").Value = "4429"
$wsInfo.Range("B2").Value = "Mitchell Josef Santner"
$wsInfo.Range("C2").Value = "Left Handed"
$wsInfo.Range("D2").Value = "Left Arm Orthodox"

# ----------------------------------------------------------------------
# 4) New "ODI Batting Extra" sheet, placed after "ODI Bowling" (last tab).
# ----------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsExtra = $wb.Worksheets.Add($null, $lastSheet)
$wsExtra.Name = "ODI Batting Extra"

$wsExtra.Range("A1").Value = "MATCH_CODE"
$wsExtra.Range("B1").Value = "BATTING_POSITION"
$wsExtra.Range("C1").Value = "NUM_4"
$wsExtra.Range("D1").Value = "NUM_6"
$wsExtra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$wsExtra.Range("F1").Value = "MAN_OF_MATCH"

$wsExtra.Range("A1:F1").Style = $wsBatting.Range("A1").Style

$extraLastRow = 21
$wsExtra.Range("A2:A$extraLastRow").NumberFormat = "@"
$wsExtra.Range("C2:E$extraLastRow").NumberFormat = "@"

$extraData = @(
    @("4453", 8, "", "", "", "NO"),
    @("4455", 8, "0", "0", "0.94%", "NO"),
    @("4608", 8, "1", "0", "2.74%", "NO"),
    @("4614", 8, "2", "0", "3.89%", "NO"),
    @("4625", 8, "", "", "", "NO"),
    @("4636", 8, "3", "0", "13.16%", "NO"),
    @("4639", 9, "1", "1", "12.26%", "NO"),
    @("4642", "", "", "", "", "NO"),
    @("4647", "", "", "", "", "NO"),
    @("4648", 8, "2", "0", "19.51%", "NO"),
    @("4649", 8, "2", "1", "12.40%", "NO"),
    @("4669", 7, "", "", "", "NO"),
    @("4673", "", "", "", "", "NO"),
    @("4676", 8, "", "", "", "NO"),
    @("4686", "", "", "", "", "NO"),
    @("4688", 8, "1", "1", "14.18%", "NO"),
    @("4690", "", "", "", "", "NO"),
    @("4692", "", "", "", "", "NO"),
    @("4695", 8, "3", "0", "25.00%", "NO"),
    @("4697", 8, "3", "2", "11.53%", "NO")
)

$r = 2
foreach ($row in $extraData) {
    $wsExtra.Cells.Item($r, 1).Value = $row[0]
    if ($row[1] -ne "") {
        $wsExtra.Cells.Item($r, 2).Value = $row[1]
    }
    $wsExtra.Cells.Item($r, 3).Value = $row[2]
    $wsExtra.Cells.Item($r, 4).Value = $row[3]
    $wsExtra.Cells.Item($r, 5).Value = $row[4]
    $wsExtra.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}
